$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variableNames")

# Insert a new row at position 4 (pushes existing rows 4-19 down to 5-20),
# and populate it with the new "month" / "F25" variable entry.
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "month"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "F25"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "F25"

# Update the view: zoom to 120% and select G4 (matches the saved view state).
$ws.Activate()
$excel.ActiveWindow.Zoom = 120
$ws.Range("G4").Select()
